$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4051.7563
$ws.Range("J51").Value = 4379.952
$ws.Range("L51").Value = 4379.952
$ws.Range("N51").Value = -5347.952
$ws.Range("H70").Value = 2042796.2
$ws.Range("I70").Value = 10204081
$ws.Range("J70").Value = 2475
$ws.Range("K70").Value = 30612243
$ws.Range("L70").Value = 7425
$ws.Range("M70").Value = -30611973
$ws.Range("N70").Value = -7965
$ws.Range("H73").Value = 2042796.2
$ws.Range("I73").Value = 10204081
$ws.Range("J73").Value = 2475
$ws.Range("K73").Value = 30612243
$ws.Range("L73").Value = 7425
$ws.Range("M73").Value = -30611307
$ws.Range("N73").Value = -9297
$ws.Range("H92").Value = 124641.25
$ws.Range("I92").Value = 411.375
$ws.Range("K92").Value = 411.375
$ws.Range("M92").Value = 836.625
$ws.Range("H98").Value = 730.2414
$ws.Range("I98").Value = 736.5
$ws.Range("K98").Value = 736.5
$ws.Range("M98").Value = 761.5
$ws.Range("H122").Value = 730.2414
$ws.Range("I122").Value = 736.5
$ws.Range("K122").Value = 2209.5
$ws.Range("M122").Value = 240.5
$ws.Range("H127").Value = 1291.55
$ws.Range("I127").Value = 1220.7858
$ws.Range("J127").Value = 1456.6666
$ws.Range("K127").Value = 3662.3574
$ws.Range("L127").Value = 4369.9998
$ws.Range("M127").Value = 1297.6426
$ws.Range("N127").Value = -14289.9998
$ws.Range("H137").Value = 27779946
$ws.Range("I137").Value = 71430240
$ws.Range("J137").Value = 2485.7727
$ws.Range("K137").Value = 214290720
$ws.Range("L137").Value = 7457.3181
$ws.Range("M137").Value = -214288170
$ws.Range("N137").Value = -12557.3181
$ws.Range("H138").Value = 1893.289
$ws.Range("I138").Value = 1611.1428
$ws.Range("J138").Value = 2358
$ws.Range("K138").Value = 4833.428400000001
$ws.Range("L138").Value = 7074
$ws.Range("M138").Value = 306.5715999999993
$ws.Range("N138").Value = -17354
$ws.Range("H141").Value = 1116.5264
$ws.Range("I141").Value = 1116.5264
$ws.Range("K141").Value = 3349.5792
$ws.Range("M141").Value = 1830.4208

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 62500
$ws.Range("I44").Value = 50000
$ws.Range("K44").Value = 50000
$ws.Range("M44").Value = -49512
$ws.Range("H51").Value = 34999
$ws.Range("I51").Value = 34999
$ws.Range("K51").Value = 34999
$ws.Range("M51").Value = -34243
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31538
$ws.Range("H132").Value = 1606.8
$ws.Range("I132").Value = 1547.0769
$ws.Range("K132").Value = 4641.2307
$ws.Range("M132").Value = -2111.2307

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4771.6113
$ws.Range("J94").Value = 2370.7144
$ws.Range("L94").Value = 2370.7144
$ws.Range("N94").Value = -3272.7144
$ws.Range("H134").Value = 1143.6154
$ws.Range("I134").Value = 1143.6154
$ws.Range("K134").Value = 3430.8462
$ws.Range("M134").Value = -895.8462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2120.1428
$ws.Range("I16").Value = 2120.1428
$ws.Range("K16").Value = 2120.1428
$ws.Range("M16").Value = -1833.1428
$ws.Range("H107").Value = 2187.9546
$ws.Range("I107").Value = 1199.3125
$ws.Range("J107").Value = 4824.3335
$ws.Range("K107").Value = 1199.3125
$ws.Range("L107").Value = 4824.3335
$ws.Range("M107").Value = 720.6875
$ws.Range("N107").Value = -8664.333500000001
$ws.Range("H113").Value = 2120.1428
$ws.Range("I113").Value = 2120.1428
$ws.Range("K113").Value = 2120.1428
$ws.Range("M113").Value = 49.85719999999992
$ws.Range("H134").Value = 3605.25
$ws.Range("I134").Value = 3263.2856
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 9789.856800000001
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = -7254.856800000001
$ws.Range("N134").Value = -23067

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 15000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -15346
$ws.Range("H33").Value = 176.2
$ws.Range("I33").Value = 176.2
$ws.Range("K33").Value = 1057.2
$ws.Range("M33").Value = -774.1999999999998
$ws.Range("H92").Value = 288.7143
$ws.Range("I92").Value = 318.75
$ws.Range("K92").Value = 956.25
$ws.Range("M92").Value = 291.75
$ws.Range("H113").Value = 1389.2858
$ws.Range("I113").Value = 774.44446
$ws.Range("J113").Value = 1850.4166
$ws.Range("K113").Value = 2323.33338
$ws.Range("L113").Value = 5551.2498
$ws.Range("M113").Value = -153.33338
$ws.Range("N113").Value = -9891.2498
$ws.Range("H116").Value = 1000
$ws.Range("I116").Value = 1000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442
$ws.Range("H131").Value = 5509.857
$ws.Range("I131").Value = 864.625
$ws.Range("J131").Value = 8368.462
$ws.Range("K131").Value = 2593.875
$ws.Range("L131").Value = 25105.386
$ws.Range("M131").Value = 2446.125
$ws.Range("N131").Value = -35185.386

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 16250
$ws.Range("J20").Value = 16250
$ws.Range("L20").Value = 16250
$ws.Range("N20").Value = -16740
$ws.Range("H24").Value = 36200
$ws.Range("J24").Value = 36200
$ws.Range("L24").Value = 36200
$ws.Range("N24").Value = -36546
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H107").Value = 813.5357
$ws.Range("I107").Value = 668.5625
$ws.Range("K107").Value = 668.5625
$ws.Range("M107").Value = 1251.4375
$ws.Range("H113").Value = 2649.8635
$ws.Range("I113").Value = 1921.8462
$ws.Range("J113").Value = 3701.4443
$ws.Range("K113").Value = 1921.8462
$ws.Range("L113").Value = 3701.4443
$ws.Range("M113").Value = 248.1538
$ws.Range("N113").Value = -8041.4443
$ws.Range("H122").Value = 1951.7273
$ws.Range("I122").Value = 1790.4546
$ws.Range("K122").Value = 5371.3638
$ws.Range("M122").Value = -2921.3638

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 733.7692
$ws.Range("I22").Value = 682.4286
$ws.Range("J22").Value = 793.6667
$ws.Range("K22").Value = 682.4286
$ws.Range("L22").Value = 793.6667
$ws.Range("M22").Value = -387.4286
$ws.Range("N22").Value = -1383.6667
$ws.Range("H27").Value = 733.7692
$ws.Range("I27").Value = 682.4286
$ws.Range("J27").Value = 793.6667
$ws.Range("K27").Value = 682.4286
$ws.Range("L27").Value = 793.6667
$ws.Range("M27").Value = -575.4286
$ws.Range("N27").Value = -1007.6667
$ws.Range("H46").Value = 2194.9119
$ws.Range("I46").Value = 1565
$ws.Range("J46").Value = 3094.7856
$ws.Range("K46").Value = 1565
$ws.Range("L46").Value = 3094.7856
$ws.Range("M46").Value = -1377
$ws.Range("N46").Value = -3470.7856
$ws.Range("H61").Value = 1700
$ws.Range("J61").Value = 1700
$ws.Range("L61").Value = 1700
$ws.Range("N61").Value = -2104
$ws.Range("H113").Value = 1700
$ws.Range("J113").Value = 1700
$ws.Range("L113").Value = 1700
$ws.Range("N113").Value = -6040
$ws.Range("H122").Value = 3007.4238
$ws.Range("I122").Value = 2276.8333
$ws.Range("J122").Value = 6195.4546
$ws.Range("K122").Value = 6830.499899999999
$ws.Range("L122").Value = 18586.3638
$ws.Range("M122").Value = -4380.499899999999
$ws.Range("N122").Value = -23486.3638

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 12999
$ws.Range("J15").Value = 12999
$ws.Range("L15").Value = 12999
$ws.Range("N15").Value = -13575
$ws.Range("H132").Value = 2734.3076
$ws.Range("J132").Value = 3640.9
$ws.Range("L132").Value = 10922.7
$ws.Range("N132").Value = -15982.7
